$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1696.6364
$ws.Range("I8").Value = 141.4
$ws.Range("J8").Value = 1974.3572
$ws.Range("K8").Value = 424.2
$ws.Range("L8").Value = 5923.071599999999
$ws.Range("M8").Value = -285.2
$ws.Range("N8").Value = -6201.071599999999
$ws.Range("H15").Value = 838.5893
$ws.Range("I15").Value = 838.5893
$ws.Range("K15").Value = 2515.7679
$ws.Range("M15").Value = -2346.7679
$ws.Range("H33").Value = 117.21429
$ws.Range("I33").Value = 119
$ws.Range("J33").Value = 116.5
$ws.Range("K33").Value = 119
$ws.Range("L33").Value = 116.5
$ws.Range("M33").Value = 110
$ws.Range("N33").Value = -574.5
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H116").Value = 11411.923
$ws.Range("J116").Value = 4420.5
$ws.Range("L116").Value = 4420.5
$ws.Range("N116").Value = -11304.5
$ws.Range("H132").Value = 1198.55
$ws.Range("I132").Value = 1109.9375
$ws.Range("K132").Value = 3329.8125
$ws.Range("M132").Value = -799.8125
$ws.Range("H137").Value = 1663.4445
$ws.Range("I137").Value = 1432.3684
$ws.Range("K137").Value = 4297.1052
$ws.Range("M137").Value = -1747.1052
$ws.Range("H138").Value = 1614.0344
$ws.Range("J138").Value = 2339.318
$ws.Range("L138").Value = 7017.954000000001
$ws.Range("N138").Value = -17297.954

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 932.63635
$ws.Range("I2").Value = 908.41174
$ws.Range("K2").Value = 908.41174
$ws.Range("M2").Value = -795.41174
$ws.Range("H17").Value = 49937.5
$ws.Range("J17").Value = 49937.5
$ws.Range("L17").Value = 49937.5
$ws.Range("N17").Value = -50283.5
$ws.Range("H74").Value = 1099.5161
$ws.Range("I74").Value = 558
$ws.Range("K74").Value = 558
$ws.Range("M74").Value = 316
$ws.Range("H77").Value = 1099.5161
$ws.Range("I77").Value = 558
$ws.Range("K77").Value = 2790
$ws.Range("M77").Value = 1578
$ws.Range("H102").Value = 1383.5
$ws.Range("I102").Value = 1161
$ws.Range("K102").Value = 1161
$ws.Range("M102").Value = 461
$ws.Range("H116").Value = 932.63635
$ws.Range("I116").Value = 908.41174
$ws.Range("K116").Value = 908.41174
$ws.Range("M116").Value = 1385.58826
$ws.Range("H132").Value = 2268.5386
$ws.Range("I132").Value = 1772.091
$ws.Range("K132").Value = 5316.272999999999
$ws.Range("M132").Value = -2786.272999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 932.63635
$ws.Range("I3").Value = 908.41174
$ws.Range("K3").Value = 908.41174
$ws.Range("M3").Value = -794.41174
$ws.Range("H107").Value = 552.5
$ws.Range("I107").Value = 484.0625
$ws.Range("K107").Value = 484.0625
$ws.Range("M107").Value = 1435.9375
$ws.Range("H134").Value = 5968.25
$ws.Range("I134").Value = 6987.8335
$ws.Range("K134").Value = 20963.5005
$ws.Range("M134").Value = -18428.5005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1627.8182
$ws.Range("I31").Value = 1181.6
$ws.Range("J31").Value = 1999.6666
$ws.Range("K31").Value = 1181.6
$ws.Range("L31").Value = 1999.6666
$ws.Range("M31").Value = -886.5999999999999
$ws.Range("N31").Value = -2589.6666
$ws.Range("H34").Value = 1627.8182
$ws.Range("I34").Value = 1181.6
$ws.Range("J34").Value = 1999.6666
$ws.Range("K34").Value = 1181.6
$ws.Range("L34").Value = 1999.6666
$ws.Range("M34").Value = -979.5999999999999
$ws.Range("N34").Value = -2403.6666
$ws.Range("H132").Value = 2237.4
$ws.Range("I132").Value = 1657.7693
$ws.Range("J132").Value = 3313.8572
$ws.Range("K132").Value = 4973.3079
$ws.Range("L132").Value = 9941.571599999999
$ws.Range("M132").Value = -2443.3079
$ws.Range("N132").Value = -15001.5716
$ws.Range("H134").Value = 1911.7059
$ws.Range("I134").Value = 1843.6875
$ws.Range("K134").Value = 5531.0625
$ws.Range("M134").Value = -2996.0625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 515.73334
$ws.Range("I7").Value = 277.6
$ws.Range("J7").Value = 634.8
$ws.Range("K7").Value = 832.8000000000001
$ws.Range("L7").Value = 1904.4
$ws.Range("M7").Value = -720.8000000000001
$ws.Range("N7").Value = -2128.4
$ws.Range("H16").Value = 3831.6667
$ws.Range("I16").Value = 3831.6667
$ws.Range("K16").Value = 11495.0001
$ws.Range("M16").Value = -11322.0001
$ws.Range("H113").Value = 6468.4443
$ws.Range("J113").Value = 964.4375
$ws.Range("L113").Value = 2893.3125
$ws.Range("N113").Value = -7233.3125
$ws.Range("H131").Value = 9889.359
$ws.Range("J131").Value = 10388.514
$ws.Range("L131").Value = 31165.542
$ws.Range("N131").Value = -41245.542
$ws.Range("H139").Value = 1834.5
$ws.Range("I139").Value = 1723.8334
$ws.Range("J139").Value = 2332.5
$ws.Range("K139").Value = 5171.5002
$ws.Range("L139").Value = 6997.5
$ws.Range("M139").Value = -31.5002000000004
$ws.Range("N139").Value = -17277.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 98.7
$ws.Range("I2").Value = 39.5
$ws.Range("K2").Value = 39.5
$ws.Range("M2").Value = 73.5
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 3973
$ws.Range("I132").Value = 3691.8
$ws.Range("J132").Value = 4324.5
$ws.Range("K132").Value = 11075.4
$ws.Range("L132").Value = 12973.5
$ws.Range("M132").Value = -8545.400000000001
$ws.Range("N132").Value = -18033.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2016.4445
$ws.Range("I82").Value = 1454
$ws.Range("J82").Value = 3985
$ws.Range("K82").Value = 1454
$ws.Range("L82").Value = 3985
$ws.Range("M82").Value = -1093
$ws.Range("N82").Value = -4707
$ws.Range("H85").Value = 2016.4445
$ws.Range("I85").Value = 1454
$ws.Range("J85").Value = 3985
$ws.Range("K85").Value = 1454
$ws.Range("L85").Value = 3985
$ws.Range("M85").Value = -206
$ws.Range("N85").Value = -6481
$ws.Range("H93").Value = 17544848
$ws.Range("I93").Value = 858.4666999999999
$ws.Range("J93").Value = 83334810
$ws.Range("K93").Value = 858.4666999999999
$ws.Range("L93").Value = 83334810
$ws.Range("M93").Value = 389.5333000000001
$ws.Range("N93").Value = -83337306
$ws.Range("H132").Value = 2051.1614
$ws.Range("I132").Value = 1559.0834
$ws.Range("J132").Value = 2361.9473
$ws.Range("K132").Value = 4677.2502
$ws.Range("L132").Value = 7085.841899999999
$ws.Range("M132").Value = -2147.2502
$ws.Range("N132").Value = -12145.8419
$ws.Range("H136").Value = 5223.579
$ws.Range("I136").Value = 4176.0713
$ws.Range("K136").Value = 12528.2139
$ws.Range("M136").Value = -9978.213899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5049.5
$ws.Range("I81").Value = 1100
$ws.Range("J81").Value = 8999
$ws.Range("K81").Value = 2200
$ws.Range("L81").Value = 17998
$ws.Range("M81").Value = -1139
$ws.Range("N81").Value = -20120
$ws.Range("H84").Value = 5049.5
$ws.Range("I84").Value = 1100
$ws.Range("J84").Value = 8999
$ws.Range("K84").Value = 11000
$ws.Range("L84").Value = 89990
$ws.Range("M84").Value = -5696
$ws.Range("N84").Value = -100598
$ws.Range("H123").Value = 47836.184
$ws.Range("J123").Value = 47836.184
$ws.Range("L123").Value = 47836.184
$ws.Range("N123").Value = -57636.184
$ws.Range("H136").Value = 1980.0646
$ws.Range("I136").Value = 1822.0869
$ws.Range("K136").Value = 5466.2607
$ws.Range("M136").Value = -2916.2607
